# Quiz Game instructions paragraph: split the single run describing how to
# answer a question into five runs (same Arial/28pt formatting) that add the
# "type a/b/c/d" and "press Enter" affordances described in the commit message.
$d = $word.ActiveDocument

$searchRange = $d.Content
$found = $searchRange.Find.Execute("Click on the button next to the answer choices to select and answer. Click on ‘NEXT’ to go to the next question.", $true, $false, $false, $false, $false, $true, 1, $false, "Click on the button", 2)
if (-not $found) {
    throw "Could not find the quiz-game instruction sentence to edit"
}

# $searchRange now spans just the replacement text ("Click on the button");
# remember where it ends so we know where to splice in the remaining runs.
$segEnd = $searchRange.End

# --- splice in the next run ---
$insertionPoint = $d.Range($segEnd, $segEnd)
[void]$insertionPoint.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p w:rsidR=`"00114789`" w:rsidRPr=`"00114789`" w:rsidRDefault=`"00114789`" w:rsidP=`"00114789`"><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:ind w:left=`"1440`"/><w:jc w:val=`"both`"/><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:t xml:space=`"preserve`"> or type a/b/c/d to </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")
# InsertXML drops the new run in as its own paragraph right after the
# insertion point; delete the paragraph mark it introduced so the new run
# rejoins the original paragraph instead of starting a new one.
$newParaMark = $d.Range($segEnd, $segEnd + 1)
$newParaMark.Delete()
$segEnd = $segEnd + 20

# --- splice in the next run ---
$insertionPoint = $d.Range($segEnd, $segEnd)
[void]$insertionPoint.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p w:rsidR=`"00114789`" w:rsidRPr=`"00114789`" w:rsidRDefault=`"00114789`" w:rsidP=`"00114789`"><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:ind w:left=`"1440`"/><w:jc w:val=`"both`"/><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:t xml:space=`"preserve`">answer. Click on ‘NEXT’ </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")
# InsertXML drops the new run in as its own paragraph right after the
# insertion point; delete the paragraph mark it introduced so the new run
# rejoins the original paragraph instead of starting a new one.
$newParaMark = $d.Range($segEnd, $segEnd + 1)
$newParaMark.Delete()
$segEnd = $segEnd + 24

# --- splice in the next run ---
$insertionPoint = $d.Range($segEnd, $segEnd)
[void]$insertionPoint.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p w:rsidR=`"00114789`" w:rsidRPr=`"00114789`" w:rsidRDefault=`"00114789`" w:rsidP=`"00114789`"><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:ind w:left=`"1440`"/><w:jc w:val=`"both`"/><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:t xml:space=`"preserve`">or press ‘Enter’ </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")
# InsertXML drops the new run in as its own paragraph right after the
# insertion point; delete the paragraph mark it introduced so the new run
# rejoins the original paragraph instead of starting a new one.
$newParaMark = $d.Range($segEnd, $segEnd + 1)
$newParaMark.Delete()
$segEnd = $segEnd + 17

# --- splice in the next run ---
$insertionPoint = $d.Range($segEnd, $segEnd)
[void]$insertionPoint.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p w:rsidR=`"00114789`" w:rsidRPr=`"00114789`" w:rsidRDefault=`"00114789`" w:rsidP=`"00114789`"><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:ind w:left=`"1440`"/><w:jc w:val=`"both`"/><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:sz w:val=`"28`"/><w:szCs w:val=`"28`"/></w:rPr><w:t>to go to the next question.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")
# InsertXML drops the new run in as its own paragraph right after the
# insertion point; delete the paragraph mark it introduced so the new run
# rejoins the original paragraph instead of starting a new one.
$newParaMark = $d.Range($segEnd, $segEnd + 1)
$newParaMark.Delete()
$segEnd = $segEnd + 27

Write-Output "Quiz instructions run split into 5 segments."
